$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$th = $nm.Theme
$cs2 = $th.ThemeColorScheme
$c = $cs2.Colors(1)
Write-Output "before $($c.RGB)"
$c.RGB = 98765
Write-Output "after $($c.RGB)"
